$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 41620
$ws.Range("D2").Value = 60091863
$ws.Range("C3").Value = 98833
$ws.Range("D3").Value = 144774505
$ws.Range("C4").Value = 33648
$ws.Range("D4").Value = 49800921
$ws.Range("C5").Value = 9621
$ws.Range("D5").Value = 14291742
$ws.Range("C6").Value = 2384
$ws.Range("D6").Value = 3541971
$ws.Range("C7").Value = 250
$ws.Range("D7").Value = 370093
$ws.Range("C12").Value = 44755
$ws.Range("D12").Value = 60538882
$ws.Range("C13").Value = 10541
$ws.Range("D13").Value = 15223418
$ws.Range("C14").Value = 27836
$ws.Range("D14").Value = 40790156
$ws.Range("C15").Value = 8822
$ws.Range("D15").Value = 13091166
$ws.Range("C16").Value = 2336
$ws.Range("D16").Value = 3469375
$ws.Range("C17").Value = 485
$ws.Range("D17").Value = 716623
$ws.Range("C20").Value = 10981
$ws.Range("D20").Value = 14453798
$ws.Range("C21").Value = 14488
$ws.Range("D21").Value = 20878818
$ws.Range("C22").Value = 33672
$ws.Range("D22").Value = 49370074
$ws.Range("C23").Value = 10844
$ws.Range("D23").Value = 16113085
$ws.Range("C24").Value = 2878
$ws.Range("D24").Value = 4277115
$ws.Range("C27").Value = 12509
$ws.Range("D27").Value = 16620862
$ws.Range("C28").Value = 8425
$ws.Range("D28").Value = 12184540
$ws.Range("C29").Value = 24235
$ws.Range("D29").Value = 35544886
$ws.Range("C30").Value = 8335
$ws.Range("D30").Value = 12391207
$ws.Range("C32").Value = 431
$ws.Range("D32").Value = 637121
$ws.Range("C34").Value = 9004
$ws.Range("D34").Value = 11858007
$ws.Range("C35").Value = 3686
$ws.Range("D35").Value = 5316568
$ws.Range("C36").Value = 8640
$ws.Range("D36").Value = 12623245
$ws.Range("C38").Value = 877
$ws.Range("D38").Value = 1306555
$ws.Range("C41").Value = 2764
$ws.Range("D41").Value = 3719772
$ws.Range("C42").Value = 18841
$ws.Range("D42").Value = 27201780
$ws.Range("C43").Value = 54833
$ws.Range("D43").Value = 80326347
$ws.Range("C44").Value = 20100
$ws.Range("D44").Value = 29835081
$ws.Range("C45").Value = 6026
$ws.Range("D45").Value = 8961336
$ws.Range("C46").Value = 1413
$ws.Range("D46").Value = 2109144
$ws.Range("C50").Value = 18283
$ws.Range("D50").Value = 24201952
$ws.Range("C52").Value = 7987
$ws.Range("D52").Value = 11732028
$ws.Range("C53").Value = 2677
$ws.Range("D53").Value = 3995633
$ws.Range("C54").Value = 845
$ws.Range("D54").Value = 1262414
$ws.Range("C57").Value = 7957
$ws.Range("D57").Value = 10947304
$ws.Range("C58").Value = 1654
$ws.Range("D58").Value = 3326816
$ws.Range("C59").Value = 3944
$ws.Range("D59").Value = 7897312
$ws.Range("C60").Value = 1560
$ws.Range("D60").Value = 3130126
$ws.Range("C61").Value = 531
$ws.Range("D61").Value = 1059083
$ws.Range("C62").Value = 193
$ws.Range("D62").Value = 397487
$ws.Range("C64").Value = 2568
$ws.Range("D64").Value = 4769464
$ws.Range("C65").Value = 16989
$ws.Range("D65").Value = 24518815
$ws.Range("C66").Value = 48348
$ws.Range("D66").Value = 70672472
$ws.Range("C67").Value = 16890
$ws.Range("D67").Value = 25096450
$ws.Range("C68").Value = 4929
$ws.Range("D68").Value = 7341024
$ws.Range("C69").Value = 1093
$ws.Range("D69").Value = 1625199
$ws.Range("C73").Value = 16136
$ws.Range("D73").Value = 21171819
$ws.Range("C74").Value = 62616
$ws.Range("D74").Value = 91009207
$ws.Range("C75").Value = 170464
$ws.Range("D75").Value = 250849987
$ws.Range("C76").Value = 72856
$ws.Range("D76").Value = 108501564
$ws.Range("C77").Value = 23850
$ws.Range("D77").Value = 35619074
$ws.Range("C78").Value = 6214
$ws.Range("D78").Value = 9273377
$ws.Range("C79").Value = 453
$ws.Range("D79").Value = 674956
$ws.Range("C85").Value = 61648
$ws.Range("D85").Value = 83113578
$ws.Range("C86").Value = 5134
$ws.Range("D86").Value = 7437373
$ws.Range("C87").Value = 12619
$ws.Range("D87").Value = 18530221
$ws.Range("C88").Value = 4137
$ws.Range("D88").Value = 6164626
$ws.Range("C89").Value = 1456
$ws.Range("D89").Value = 2175111
$ws.Range("C90").Value = 350
$ws.Range("D90").Value = 521512
$ws.Range("C91").Value = 33
$ws.Range("D91").Value = 49402
$ws.Range("C93").Value = 5879
$ws.Range("D93").Value = 7885374
$ws.Range("C94").Value = 1851
$ws.Range("D94").Value = 2664993
$ws.Range("C95").Value = 5876
$ws.Range("D95").Value = 8658589
$ws.Range("C101").Value = 3989
$ws.Range("D101").Value = 5286563
$ws.Range("C107").Value = 11907
$ws.Range("D107").Value = 17255109
$ws.Range("C108").Value = 31134
$ws.Range("D108").Value = 45694469
$ws.Range("C109").Value = 10436
$ws.Range("D109").Value = 15513260
$ws.Range("C110").Value = 2895
$ws.Range("D110").Value = 4315571
$ws.Range("C112").Value = 68
$ws.Range("D112").Value = 102000
$ws.Range("C115").Value = 10497
$ws.Range("D115").Value = 13808872
$ws.Range("C116").Value = 33447
$ws.Range("D116").Value = 48181144
$ws.Range("C117").Value = 70957
$ws.Range("D117").Value = 103773661
$ws.Range("C118").Value = 22751
$ws.Range("D118").Value = 33791638
$ws.Range("C119").Value = 6532
$ws.Range("D119").Value = 9724052
$ws.Range("C120").Value = 1307
$ws.Range("D120").Value = 1953233
$ws.Range("C121").Value = 130
$ws.Range("D121").Value = 189795
$ws.Range("C125").Value = 27716
$ws.Range("D125").Value = 36896376
$ws.Range("C126").Value = 39955
$ws.Range("D126").Value = 57594589
$ws.Range("C127").Value = 83393
$ws.Range("D127").Value = 121830127
$ws.Range("C128").Value = 25547
$ws.Range("D128").Value = 37904937
$ws.Range("C129").Value = 6969
$ws.Range("D129").Value = 10356009
$ws.Range("C130").Value = 1479
$ws.Range("D130").Value = 2191761
$ws.Range("C134").Value = 34231
$ws.Range("D134").Value = 45316663
$ws.Range("C135").Value = 14524
$ws.Range("D135").Value = 21007024
$ws.Range("C136").Value = 34634
$ws.Range("D136").Value = 50827850
$ws.Range("C137").Value = 12216
$ws.Range("D137").Value = 18150271
$ws.Range("C138").Value = 3233
$ws.Range("D138").Value = 4819375
$ws.Range("C139").Value = 589
$ws.Range("D139").Value = 877490
$ws.Range("C143").Value = 11608
$ws.Range("D143").Value = 15420764
$ws.Range("C144").Value = 39250
$ws.Range("D144").Value = 56662334
$ws.Range("C145").Value = 89880
$ws.Range("D145").Value = 131572025
$ws.Range("C146").Value = 26873
$ws.Range("D146").Value = 39923880
$ws.Range("C147").Value = 7165
$ws.Range("D147").Value = 10677931
$ws.Range("C148").Value = 1762
$ws.Range("D148").Value = 2623218
$ws.Range("C149").Value = 120
$ws.Range("D149").Value = 179630
$ws.Range("C151").Value = 31738
$ws.Range("D151").Value = 42646669
